$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 985.5
$ws.Range("J17").Value = 985.5
$ws.Range("L17").Value = 2956.5
$ws.Range("N17").Value = -3292.5
# Row 62
$ws.Range("H62").Value = 5058.7036
$ws.Range("I62").Value = 4636.7
$ws.Range("K62").Value = 4636.7
$ws.Range("M62").Value = -4012.7
# Row 65
$ws.Range("H65").Value = 5058.7036
$ws.Range("I65").Value = 4636.7
$ws.Range("K65").Value = 23183.5
$ws.Range("M65").Value = -20063.5
# Row 76
$ws.Range("H76").Value = 4823.4287
$ws.Range("I76").Value = 4216.6
$ws.Range("K76").Value = 4216.6
$ws.Range("M76").Value = -3901.6
# Row 79
$ws.Range("H79").Value = 4823.4287
$ws.Range("I79").Value = 4216.6
$ws.Range("K79").Value = 4216.6
$ws.Range("M79").Value = -3124.6
# Row 88
$ws.Range("H88").Value = 15248.5
$ws.Range("I88").Value = 3500
$ws.Range("K88").Value = 3500
$ws.Range("M88").Value = -3094
# Row 91
$ws.Range("H91").Value = 15248.5
$ws.Range("I91").Value = 3500
$ws.Range("K91").Value = 3500
$ws.Range("M91").Value = -2096
# Row 113
$ws.Range("H113").Value = 8956.286
$ws.Range("J113").Value = 9600
$ws.Range("L113").Value = 9600
$ws.Range("N113").Value = -16108
# Row 125
$ws.Range("H125").Value = 26609.691
$ws.Range("I125").Value = 53071.832
$ws.Range("J125").Value = 3927.8572
$ws.Range("K125").Value = 477646.488
$ws.Range("L125").Value = 35350.7148
$ws.Range("M125").Value = -475186.488
$ws.Range("N125").Value = -40270.7148

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3781.5715
$ws.Range("I45").Value = 1343.2858
$ws.Range("K45").Value = 1343.2858
$ws.Range("M45").Value = -966.2858000000001
# Row 63
$ws.Range("H63").Value = 2251.1428
$ws.Range("I63").Value = 1893
$ws.Range("K63").Value = 1893
$ws.Range("M63").Value = -1207
# Row 66
$ws.Range("H66").Value = 2251.1428
$ws.Range("I66").Value = 1893
$ws.Range("K66").Value = 9465
$ws.Range("M66").Value = -6033
# Row 74
$ws.Range("H74").Value = 278773.47
$ws.Range("I74").Value = 462215.78
$ws.Range("K74").Value = 462215.78
$ws.Range("M74").Value = -461341.78
# Row 77
$ws.Range("H77").Value = 278773.47
$ws.Range("I77").Value = 462215.78
$ws.Range("K77").Value = 2311078.9
$ws.Range("M77").Value = -2306710.9
# Row 102
$ws.Range("H102").Value = 4978.778
$ws.Range("I102").Value = 4912.25
$ws.Range("K102").Value = 4912.25
$ws.Range("M102").Value = -3290.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3155.375
$ws.Range("I20").Value = 3879.7144
$ws.Range("J20").Value = 1465.25
$ws.Range("K20").Value = 3879.7144
$ws.Range("L20").Value = 1465.25
$ws.Range("M20").Value = -3632.7144
$ws.Range("N20").Value = -1959.25
# Row 86
$ws.Range("H86").Value = 1987.1875
$ws.Range("I86").Value = 1747.25
$ws.Range("J86").Value = 2067.1667
$ws.Range("K86").Value = 1747.25
$ws.Range("L86").Value = 2067.1667
$ws.Range("M86").Value = -624.25
$ws.Range("N86").Value = -4313.1667
# Row 89
$ws.Range("H89").Value = 1987.1875
$ws.Range("I89").Value = 1747.25
$ws.Range("J89").Value = 2067.1667
$ws.Range("K89").Value = 8736.25
$ws.Range("L89").Value = 10335.8335
$ws.Range("M89").Value = -3120.25
$ws.Range("N89").Value = -21567.8335
# Row 107
$ws.Range("H107").Value = 3711
$ws.Range("I107").Value = 4079.5293
$ws.Range("J107").Value = 3229.077
$ws.Range("K107").Value = 4079.5293
$ws.Range("L107").Value = 3229.077
$ws.Range("M107").Value = -2159.5293
$ws.Range("N107").Value = -7069.077
# Row 134
$ws.Range("H134").Value = 3088.125
$ws.Range("I134").Value = 2569.8
$ws.Range("J134").Value = 3952
$ws.Range("K134").Value = 7709.400000000001
$ws.Range("L134").Value = 11856
$ws.Range("M134").Value = -5174.400000000001
$ws.Range("N134").Value = -16926

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 880
$ws.Range("J22").Value = 1050
$ws.Range("L22").Value = 1050
$ws.Range("N22").Value = -1750
# Row 31
$ws.Range("H31").Value = 5264481
$ws.Range("I31").Value = 5883244
$ws.Range("K31").Value = 5883244
$ws.Range("M31").Value = -5882949
# Row 34
$ws.Range("H34").Value = 5264481
$ws.Range("I34").Value = 5883244
$ws.Range("K34").Value = 5883244
$ws.Range("M34").Value = -5883042
# Row 62
$ws.Range("H62").Value = 5720.6875
$ws.Range("I62").Value = 3978.1
$ws.Range("J62").Value = 8625
$ws.Range("K62").Value = 3978.1
$ws.Range("L62").Value = 8625
$ws.Range("M62").Value = -3354.1
$ws.Range("N62").Value = -9873
# Row 65
$ws.Range("H65").Value = 5720.6875
$ws.Range("I65").Value = 3978.1
$ws.Range("J65").Value = 8625
$ws.Range("K65").Value = 19890.5
$ws.Range("L65").Value = 43125
$ws.Range("M65").Value = -16770.5
$ws.Range("N65").Value = -49365
# Row 107
$ws.Range("H107").Value = 1355.4445
$ws.Range("I107").Value = 1275
$ws.Range("J107").Value = 1456
$ws.Range("K107").Value = 1275
$ws.Range("L107").Value = 1456
$ws.Range("M107").Value = 645
$ws.Range("N107").Value = -5296

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 673.3333
$ws.Range("I8").Value = 673.3333
$ws.Range("K8").Value = 2019.9999
$ws.Range("M8").Value = -1880.9999
# Row 12
$ws.Range("H12").Value = 219.66667
$ws.Range("J12").Value = 210.76923
$ws.Range("L12").Value = 632.30769
$ws.Range("N12").Value = -978.30769
# Row 23
$ws.Range("H23").Value = 342.6154
$ws.Range("J23").Value = 431.125
$ws.Range("L23").Value = 1293.375
$ws.Range("N23").Value = -1763.375
# Row 117
$ws.Range("H117").Value = 1215.1428
$ws.Range("I117").Value = 1001.4
$ws.Range("J117").Value = 1749.5
$ws.Range("K117").Value = 3004.2
$ws.Range("L117").Value = 5248.5
$ws.Range("M117").Value = 437.8000000000002
$ws.Range("N117").Value = -12132.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8296.799999999999
$ws.Range("I70").Value = 9577.200000000001
$ws.Range("K70").Value = 9577.200000000001
$ws.Range("M70").Value = -9307.200000000001
# Row 73
$ws.Range("H73").Value = 8296.799999999999
$ws.Range("I73").Value = 9577.200000000001
$ws.Range("K73").Value = 9577.200000000001
$ws.Range("M73").Value = -8641.200000000001
# Row 97
$ws.Range("H97").Value = 1398.85
$ws.Range("I97").Value = 1399.3334
$ws.Range("K97").Value = 1399.3334
$ws.Range("M97").Value = -903.3334
# Row 107
$ws.Range("H107").Value = 710.7778
$ws.Range("J107").Value = 787.25
$ws.Range("L107").Value = 787.25
$ws.Range("N107").Value = -4627.25
# Row 113
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170
# Row 126
$ws.Range("H126").Value = 2712
$ws.Range("I126").Value = 1730.6666
$ws.Range("K126").Value = 5191.9998
$ws.Range("M126").Value = -2721.9998
# Row 130
$ws.Range("H130").Value = 112499.5
$ws.Range("J130").Value = 112499.5
$ws.Range("L130").Value = 112499.5
$ws.Range("N130").Value = -122539.5
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 3373333
$ws.Range("J2").Value = 3373333
$ws.Range("L2").Value = 3373333
$ws.Range("N2").Value = -3373557
# Row 7
$ws.Range("H7").Value = 3284.7144
$ws.Range("I7").Value = 3284.7144
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3284.7144
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3172.7144
$ws.Range("N7").ClearContents()
# Row 22
$ws.Range("I22").Value = 1713.5714
$ws.Range("J22").Value = 1900
$ws.Range("K22").Value = 1713.5714
$ws.Range("L22").Value = 1900
$ws.Range("M22").Value = -1418.5714
$ws.Range("N22").Value = -2490
# Row 27
$ws.Range("I27").Value = 1713.5714
$ws.Range("J27").Value = 1900
$ws.Range("K27").Value = 1713.5714
$ws.Range("L27").Value = 1900
$ws.Range("M27").Value = -1606.5714
$ws.Range("N27").Value = -2114
# Row 40
$ws.Range("H40").Value = 3000.4
$ws.Range("I40").Value = 2999.25
$ws.Range("J40").Value = 3005
$ws.Range("K40").Value = 2999.25
$ws.Range("L40").Value = 3005
$ws.Range("M40").Value = -2863.25
$ws.Range("N40").Value = -3277
# Row 55
$ws.Range("H55").Value = 3642.5
$ws.Range("I55").Value = 3150.111
$ws.Range("J55").Value = 4045.3635
$ws.Range("K55").Value = 3150.111
$ws.Range("L55").Value = 4045.3635
$ws.Range("M55").Value = -2977.111
$ws.Range("N55").Value = -4391.363499999999
# Row 93
$ws.Range("H93").Value = 2406.7144
$ws.Range("I93").Value = 1974.6666
$ws.Range("J93").Value = 4999
$ws.Range("K93").Value = 1974.6666
$ws.Range("L93").Value = 4999
$ws.Range("M93").Value = -726.6666
$ws.Range("N93").Value = -7495
# Row 126
$ws.Range("H126").Value = 3284.7144
$ws.Range("I126").Value = 3284.7144
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9854.143199999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7384.143199999999
$ws.Range("N126").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 737.63635
$ws.Range("I100").Value = 843.4286
$ws.Range("K100").Value = 1686.8572
$ws.Range("M100").Value = -1145.8572
# Row 107
$ws.Range("H107").Value = 963.3570999999999
$ws.Range("I107").Value = 932.5
$ws.Range("K107").Value = 2797.5
$ws.Range("M107").Value = -877.5
# Row 122
$ws.Range("H122").Value = 35465.8
$ws.Range("I122").Value = 39060.64
$ws.Range("J122").Value = 3112.25
$ws.Range("K122").Value = 117181.92
$ws.Range("L122").Value = 9336.75
$ws.Range("M122").Value = -114731.92
$ws.Range("N122").Value = -14236.75
# Row 125
$ws.Range("H125").Value = 50753.77
$ws.Range("J125").Value = 50753.77
$ws.Range("L125").Value = 50753.77
$ws.Range("N125").Value = -60593.77
